$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.0653562112283128
$ws.Range("C2").Value = 0.5162590491173051
$ws.Range("D2").Value = -0.2635820013437095
$ws.Range("E2").Value = 0.1328695788870839

$ws.Range("B3").Value = 0.01725040393772195
$ws.Range("C3").Value = 0.8934473402079195
$ws.Range("D3").Value = -0.2364422918969799
$ws.Range("E3").Value = 0.2709430997724238

$ws.Range("B4").Value = 0.01564798175618928
$ws.Range("C4").Value = 0.9032937125331343
$ws.Range("D4").Value = -0.2380447140785127
$ws.Range("E4").Value = 0.2693406775908912

$ws.Range("B5").Value = 0.009818463971774629
$ws.Range("C5").Value = 0.9392296143246229
$ws.Range("D5").Value = -0.2438742318629273
$ws.Range("E5").Value = 0.2635111598064765

$ws.Range("B6").Value = 0.01161460772199954
$ws.Range("C6").Value = 0.9281405609428504
$ws.Range("D6").Value = -0.2420780881127023
$ws.Range("E6").Value = 0.2653073035567014
